$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 20003680
$ws.Range("I94").Value = 33335468
$ws.Range("K94").Value = 33335468
$ws.Range("M94").Value = -33335017
$ws.Range("H97").Value = 2345
$ws.Range("J97").Value = 2345
$ws.Range("L97").Value = 7035
$ws.Range("N97").Value = -8027
$ws.Range("H103").Value = 45454956
$ws.Range("I103").Value = 401
$ws.Range("J103").Value = 71428984
$ws.Range("K103").Value = 1203
$ws.Range("L103").Value = 214286952
$ws.Range("M103").Value = -617
$ws.Range("N103").Value = -214288124
$ws.Range("H106").Value = 5478.077
$ws.Range("I106").Value = 1321.7222
$ws.Range("J106").Value = 14829.875
$ws.Range("K106").Value = 1321.7222
$ws.Range("L106").Value = 14829.875
$ws.Range("M106").Value = -690.7221999999999
$ws.Range("N106").Value = -16091.875
$ws.Range("H116").Value = 4663.5
$ws.Range("I116").Value = 4461.6
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4461.6
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1019.6
$ws.Range("N116").Value = -11884
$ws.Range("H127").Value = 700000
$ws.Range("I127").Value = 700000
$ws.Range("K127").Value = 2100000
$ws.Range("M127").Value = -2095040
$ws.Range("H134").Value = 124748
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 1967.5625
$ws.Range("I137").Value = 1581.0834
$ws.Range("J137").Value = 3127
$ws.Range("K137").Value = 4743.2502
$ws.Range("L137").Value = 9381
$ws.Range("M137").Value = -2193.2502
$ws.Range("N137").Value = -14481
$ws.Range("H138").Value = 2472.2856
$ws.Range("J138").Value = 3553.6875
$ws.Range("L138").Value = 10661.0625
$ws.Range("N138").Value = -20941.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16607.723
$ws.Range("I32").Value = 2997.3943
$ws.Range("K32").Value = 2997.3943
$ws.Range("M32").Value = -2710.3943
$ws.Range("H45").Value = 563643.9399999999
$ws.Range("I45").Value = 778960.7
$ws.Range("J45").Value = 3820.4
$ws.Range("K45").Value = 778960.7
$ws.Range("L45").Value = 3820.4
$ws.Range("M45").Value = -778583.7
$ws.Range("N45").Value = -4574.4
$ws.Range("H61").Value = 2596.6829
$ws.Range("I61").Value = 2368.1562
$ws.Range("J61").Value = 3409.2222
$ws.Range("K61").Value = 2368.1562
$ws.Range("L61").Value = 3409.2222
$ws.Range("M61").Value = -2156.1562
$ws.Range("N61").Value = -3833.2222
$ws.Range("H62").Value = 24000
$ws.Range("J62").Value = 24000
$ws.Range("L62").Value = 24000
$ws.Range("N62").Value = -25248
$ws.Range("H63").Value = 2662
$ws.Range("I63").Value = 2662
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2662
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1976
$ws.Range("N63").ClearContents()
$ws.Range("H65").Value = 24000
$ws.Range("J65").Value = 24000
$ws.Range("L65").Value = 72000
$ws.Range("N65").Value = -78240
$ws.Range("H66").Value = 2662
$ws.Range("I66").Value = 2662
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 13310
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -9878
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 2177.652
$ws.Range("I122").Value = 2214.35
$ws.Range("K122").Value = 6643.049999999999
$ws.Range("M122").Value = -4193.049999999999
$ws.Range("H134").Value = 69998
$ws.Range("J134").Value = 69998
$ws.Range("L134").Value = 69998
$ws.Range("N134").Value = -80138
$ws.Range("H136").Value = 2596.6829
$ws.Range("I136").Value = 2368.1562
$ws.Range("J136").Value = 3409.2222
$ws.Range("K136").Value = 7104.4686
$ws.Range("L136").Value = 10227.6666
$ws.Range("M136").Value = -4554.4686
$ws.Range("N136").Value = -15327.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19620
$ws.Range("J82").Value = 24998.334
$ws.Range("L82").Value = 24998.334
$ws.Range("N82").Value = -25764.334
$ws.Range("H85").Value = 19620
$ws.Range("J85").Value = 24998.334
$ws.Range("L85").Value = 24998.334
$ws.Range("N85").Value = -27650.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3016508.5
$ws.Range("I6").Value = 3519085
$ws.Range("K6").Value = 3519085
$ws.Range("M6").Value = -3518972
$ws.Range("H7").Value = 30303266
$ws.Range("I7").Value = 37037230
$ws.Range("K7").Value = 37037230
$ws.Range("M7").Value = -37037117
$ws.Range("H20").Value = 100780
$ws.Range("J20").Value = 100780
$ws.Range("L20").Value = 100780
$ws.Range("N20").Value = -101252
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H30").Value = 100780
$ws.Range("J30").Value = 100780
$ws.Range("L30").Value = 100780
$ws.Range("N30").Value = -100962
$ws.Range("H31").Value = 2851.5264
$ws.Range("I31").Value = 2719.7
$ws.Range("K31").Value = 2719.7
$ws.Range("M31").Value = -2424.7
$ws.Range("H34").Value = 2851.5264
$ws.Range("I34").Value = 2719.7
$ws.Range("K34").Value = 2719.7
$ws.Range("M34").Value = -2517.7
$ws.Range("H74").Value = 46002.1
$ws.Range("J74").Value = 49099.875
$ws.Range("L74").Value = 49099.875
$ws.Range("N74").Value = -50847.875
$ws.Range("H77").Value = 46002.1
$ws.Range("J77").Value = 49099.875
$ws.Range("L77").Value = 147299.625
$ws.Range("N77").Value = -156035.625
$ws.Range("H105").Value = 914.3333
$ws.Range("I105").Value = 860.2632
$ws.Range("J105").Value = 1119.8
$ws.Range("K105").Value = 860.2632
$ws.Range("L105").Value = 1119.8
$ws.Range("M105").Value = 886.7368
$ws.Range("N105").Value = -4613.8
$ws.Range("H128").Value = 100780
$ws.Range("J128").Value = 100780
$ws.Range("L128").Value = 100780
$ws.Range("N128").Value = -110740
$ws.Range("H133").Value = 69338.11
$ws.Range("J133").Value = 69338.11
$ws.Range("L133").Value = 69338.11
$ws.Range("N133").Value = -74398.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 729.9167
$ws.Range("I5").Value = 371.4
$ws.Range("J5").Value = 986
$ws.Range("K5").Value = 1114.2
$ws.Range("L5").Value = 2958
$ws.Range("M5").Value = -1002.2
$ws.Range("N5").Value = -3182
$ws.Range("H37").Value = 100045900
$ws.Range("J37").Value = 100045900
$ws.Range("L37").Value = 300137700
$ws.Range("N37").Value = -300137924
$ws.Range("H116").Value = 8666.125
$ws.Range("I116").Value = 8475.786
$ws.Range("J116").Value = 9998.5
$ws.Range("K116").Value = 25427.358
$ws.Range("L116").Value = 29995.5
$ws.Range("M116").Value = -21985.358
$ws.Range("N116").Value = -36879.5
$ws.Range("H135").Value = 729.9167
$ws.Range("I135").Value = 371.4
$ws.Range("J135").Value = 986
$ws.Range("K135").Value = 3342.6
$ws.Range("L135").Value = 8874
$ws.Range("M135").Value = -807.5999999999999
$ws.Range("N135").Value = -13944

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3007.4443
$ws.Range("I102").Value = 2959.5518
$ws.Range("K102").Value = 2959.5518
$ws.Range("M102").Value = -1337.5518

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9346.1875
$ws.Range("I122").Value = 9967.143
$ws.Range("K122").Value = 29901.429
$ws.Range("M122").Value = -27451.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 26322324
$ws.Range("I107").Value = 9336.75
$ws.Range("K107").Value = 28010.25
$ws.Range("M107").Value = -26090.25
$ws.Range("H122").Value = 1795.8975
$ws.Range("I122").Value = 1702.8857
$ws.Range("K122").Value = 5108.6571
$ws.Range("M122").Value = -2658.6571
$ws.Range("H132").Value = 4105.1304
$ws.Range("I132").Value = 4423.675
$ws.Range("K132").Value = 13271.025
$ws.Range("M132").Value = -10741.025
